$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "236.658.7.88/9"
$ws.Range("B3").Value = "Célio Vetrano"
$ws.Range("C3").Value = "20E"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "01/05/2019"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "{`"segunda-feira`": `"00:00 \u00e0s 00:00`", `"ter\u00e7a-feira`": `"00:00 \u00e0s 00:00`", `"quarta-feira`": `"00:00 \u00e0s 00:00`", `"quinta-feira`": `"00:00 \u00e0s 00:00`", `"sexta-feira`": `"00:00 \u00e0s 00:00`"}"
$ws.Range("F3").Value = "{`"segunda-feira`": `"12:00 \u00e0s 15:00`", `"ter\u00e7a-feira`": `"12:00 \u00e0s 15:00`", `"quarta-feira`": `"12:00 \u00e0s 15:00`", `"quinta-feira`": `"12:00 \u00e0s 15:00`", `"sexta-feira`": `"12:00 \u00e0s 15:00`"}"
$ws.Range("G3").Value = "7E"
$ws.Range("H3").Value = "Professor de Informática"

$ws.Range("A4").Value = "777.777.7.77/7"
$ws.Range("B4").Value = "Aline Silva"
$ws.Range("C4").Value = "15A"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "01/06/2018"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "{`"segunda-feira`": `"00:00 \u00e0s 00:00`", `"ter\u00e7a-feira`": `"00:00 \u00e0s 00:00`", `"quarta-feira`": `"00:00 \u00e0s 00:00`", `"quinta-feira`": `"00:00 \u00e0s 00:00`", `"sexta-feira`": `"00:00 \u00e0s 00:00`"}"
$ws.Range("F4").Value = "{`"segunda-feira`": `"08:00 \u00e0s 12:00`", `"ter\u00e7a-feira`": `"08:00 \u00e0s 12:00`", `"quarta-feira`": `"08:00 \u00e0s 12:00`", `"quinta-feira`": `"08:00 \u00e0s 12:00`", `"sexta-feira`": `"08:00 \u00e0s 12:00`"}"
$ws.Range("G4").Value = "7B"
$ws.Range("H4").Value = "professor de educação física"

$ws.Range("A5").Value = "666.666.6.66/6"
$ws.Range("B5").Value = "Alice Burba"
$ws.Range("C5").Value = "20E"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "25/06/1985"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "{`"segunda-feira`": `"13:00 \u00e0s 14:00`", `"ter\u00e7a-feira`": `"13:00 \u00e0s 14:00`", `"quarta-feira`": `"13:00 \u00e0s 14:00`", `"quinta-feira`": `"13:00 \u00e0s 14:00`", `"sexta-feira`": `"13:00 \u00e0s 14:00`"}"
$ws.Range("F5").Value = "{`"segunda-feira`": `"08:00 \u00e0s 12:00`", `"ter\u00e7a-feira`": `"08:00 \u00e0s 12:00`", `"quarta-feira`": `"08:00 \u00e0s 12:00`", `"quinta-feira`": `"08:00 \u00e0s 12:00`", `"sexta-feira`": `"08:00 \u00e0s 14:00`"}"
$ws.Range("G5").Value = "6D"
$ws.Range("H5").Value = "Professor de Educação Infantil"

$ws.Range("A6").Value = "555.555.5.55/5"
$ws.Range("B6").Value = "Tatiana Marques"
$ws.Range("C6").Value = "14B"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "01/01/2025"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "{`"segunda-feira`": `"14:00 \u00e0s 15:00`", `"ter\u00e7a-feira`": `"14:00 \u00e0s 15:00`", `"quarta-feira`": `"14:00 \u00e0s 15:00`", `"quinta-feira`": `"14:00 \u00e0s 15:00`", `"sexta-feira`": `"14:00 \u00e0s 15:00`"}"
$ws.Range("F6").Value = "{`"segunda-feira`": `"14:00 \u00e0s 15:00`", `"ter\u00e7a-feira`": `"14:00 \u00e0s 15:00`", `"quarta-feira`": `"14:00 \u00e0s 15:00`", `"quinta-feira`": `"14:00 \u00e0s 15:00`", `"sexta-feira`": `"14:00 \u00e0s 15:00`"}"
$ws.Range("G6").Value = "7A"
$ws.Range("H6").Value = "professor"

$ws.Range("A7").Value = "333.444.5.66/7"
$ws.Range("B7").Value = "Marilucia Junqueira"
$ws.Range("C7").Value = "20D"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "03/07/1996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "{`"segunda-feira`": `"00:00 \u00e0s 00:00`", `"ter\u00e7a-feira`": `"00:00 \u00e0s 00:00`", `"quarta-feira`": `"00:00 \u00e0s 00:00`", `"quinta-feira`": `"00:00 \u00e0s 00:00`", `"sexta-feira`": `"00:00 \u00e0s 00:00`"}"
$ws.Range("F7").Value = "{`"segunda-feira`": `"07:00 \u00e0s 11:00`", `"ter\u00e7a-feira`": `"07:00 \u00e0s 11:00`", `"quarta-feira`": `"07:00 \u00e0s 11:00`", `"quinta-feira`": `"07:00 \u00e0s 11:00`", `"sexta-feira`": `"07:00 \u00e0s 11:00`"}"
$ws.Range("G7").Value = "9A"
$ws.Range("H7").Value = "professor"

$ws.Range("A8").Value = "999.888.7.66/5"
$ws.Range("B8").Value = "João Paulo Sarmento"
$ws.Range("C8").Value = "22A"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "05/04/1987"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "{`"segunda-feira`": `"09:00 \u00e0s 10:00`", `"ter\u00e7a-feira`": `"09:00 \u00e0s 10:00`", `"quarta-feira`": `"09:00 \u00e0s 10:00`", `"quinta-feira`": `"09:00 \u00e0s 10:00`", `"sexta-feira`": `"09:00 \u00e0s 10:00`"}"
$ws.Range("F8").Value = "{`"segunda-feira`": `"11:00 \u00e0s 15:00`", `"ter\u00e7a-feira`": `"10:00 \u00e0s 15:00`", `"quarta-feira`": `"10:00 \u00e0s 15:00`", `"quinta-feira`": `"10:00 \u00e0s 15:00`", `"sexta-feira`": `"10:00 \u00e0s 15:00`"}"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "00"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = "Professor Readaptado"
